# Update the "LoginDetails" sheet with new test data (more details / rows of
# sample credentials), matching the commit "Project is updated with some more
# details".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginDetails")

# Font color used throughout this sheet for e-mail / credential cells
# (matches the existing "FF2A00FF" Consolas font already used in the file).
$credentialColor = 16711722   # RGB(42, 0, 255) -> 0x2A + 0x00*256 + 0xFF*65536

# --- Row 2 : Username / Password pair (plain, non-hyperlinked) -------------
$ws.Range("A2").Value = "uaaxv@gmail.com"
$ws.Range("A2").Font.Color = $credentialColor
$ws.Range("A2").Font.Underline = $true

$ws.Range("B2").Value = "LWW@Ujt"
$ws.Range("B2").Font.Color = $credentialColor

# --- Row 3 : Username / Password pair (password cell is a hyperlink) -------
$ws.Range("A3").Value = "uaaxv@gmail.com"
$ws.Range("A3").Font.Color = $credentialColor
$ws.Range("A3").Font.Underline = $true

$ws.Range("B3").Value = "LWW@Uj"
$ws.Range("B3").Style = "Hyperlink"

# --- Row 4 : Username / Password pair (plain, non-hyperlinked) -------------
$ws.Range("A4").Value = "uaaxv@gmail.com"
$ws.Range("A4").Font.Color = $credentialColor
$ws.Range("A4").Font.Underline = $true

$ws.Range("B4").Value = "LWW@Ujt"
$ws.Range("B4").Font.Color = $credentialColor

# --- Row 5 : Username / Password pair (password cell is a hyperlink) -------
$ws.Range("A5").Value = "uaaxv@gmail.com"
$ws.Range("A5").Font.Color = $credentialColor
$ws.Range("A5").Font.Underline = $true

$ws.Range("B5").Value = "LWW@Uj"
$ws.Range("B5").Style = "Hyperlink"

# Auto-fit the columns now that the data has changed width, then leave the
# selection on C3 as the workbook was last saved.
$ws.Columns("A:C").AutoFit()

$ws.Activate()
$ws.Range("C3").Select()
